{"js": "// Update the date line and the twenty-five \"N\u00d7N=NNNN\" multiplication\n// answers in the table to the new values from the target revision.\nconst replacements = [\n  [\"2024-10-03 Thursday\", \"2024-10-04 Friday\"],\n  [\"883\u00d75=4415\", \"527\u00d76=3162\"],\n  [\"195\u00d78=1560\", \"485\u00d79=4365\"],\n  [\"883\u00d73=2649\", \"658\u00d79=5922\"],\n  [\"435\u00d79=3915\", \"625\u00d79=5625\"],\n  [\"525\u00d72=1050\", \"586\u00d74=2344\"],\n  [\"961\u00d75=4805\", \"380\u00d75=1900\"],\n  [\"767\u00d73=2301\", \"509\u00d77=3563\"],\n  [\"872\u00d78=6976\", \"614\u00d76=3684\"],\n  [\"729\u00d78=5832\", \"143\u00d75=715\"],\n  [\"251\u00d74=1004\", \"970\u00d73=2910\"],\n  [\"519\u00d74=2076\", \"500\u00d77=3500\"],\n  [\"547\u00d77=3829\", \"264\u00d76=1584\"],\n  [\"811\u00d74=3244\", \"449\u00d75=2245\"],\n  [\"516\u00d76=3096\", \"921\u00d73=2763\"],\n  [\"517\u00d74=2068\", \"118\u00d77=826\"],\n  [\"519\u00d76=3114\", \"295\u00d72=590\"],\n  [\"563\u00d74=2252\", \"855\u00d76=5130\"],\n  [\"621\u00d79=5589\", \"948\u00d76=5688\"],\n  [\"101\u00d78=808\", \"456\u00d73=1368\"],\n  [\"224\u00d79=2016\", \"929\u00d77=6503\"],\n  [\"722\u00d74=2888\", \"584\u00d78=4672\"],\n  [\"228\u00d73=684\", \"164\u00d72=328\"],\n  [\"719\u00d75=3595\", \"393\u00d78=3144\"],\n  [\"122\u00d73=366\", \"258\u00d78=2064\"],\n  [\"715\u00d73=2145\", \"676\u00d76=4056\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the twenty-five \"N\u00d7N=NNNN\" multiplication\n# answers in the table to the new values from the target revision.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-03 Thursday\", \"2024-10-04 Friday\"),\n    @(\"883\u00d75=4415\", \"527\u00d76=3162\"),\n    @(\"195\u00d78=1560\", \"485\u00d79=4365\"),\n    @(\"883\u00d73=2649\", \"658\u00d79=5922\"),\n    @(\"435\u00d79=3915\", \"625\u00d79=5625\"),\n    @(\"525\u00d72=1050\", \"586\u00d74=2344\"),\n    @(\"961\u00d75=4805\", \"380\u00d75=1900\"),\n    @(\"767\u00d73=2301\", \"509\u00d77=3563\"),\n    @(\"872\u00d78=6976\", \"614\u00d76=3684\"),\n    @(\"729\u00d78=5832\", \"143\u00d75=715\"),\n    @(\"251\u00d74=1004\", \"970\u00d73=2910\"),\n    @(\"519\u00d74=2076\", \"500\u00d77=3500\"),\n    @(\"547\u00d77=3829\", \"264\u00d76=1584\"),\n    @(\"811\u00d74=3244\", \"449\u00d75=2245\"),\n    @(\"516\u00d76=3096\", \"921\u00d73=2763\"),\n    @(\"517\u00d74=2068\", \"118\u00d77=826\"),\n    @(\"519\u00d76=3114\", \"295\u00d72=590\"),\n    @(\"563\u00d74=2252\", \"855\u00d76=5130\"),\n    @(\"621\u00d79=5589\", \"948\u00d76=5688\"),\n    @(\"101\u00d78=808\", \"456\u00d73=1368\"),\n    @(\"224\u00d79=2016\", \"929\u00d77=6503\"),\n    @(\"722\u00d74=2888\", \"584\u00d78=4672\"),\n    @(\"228\u00d73=684\", \"164\u00d72=328\"),\n    @(\"719\u00d75=3595\", \"393\u00d78=3144\"),\n    @(\"122\u00d73=366\", \"258\u00d78=2064\"),\n    @(\"715\u00d73=2145\", \"676\u00d76=4056\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
